# Add "cfop" columns/values to the "PI hours" sheet and create a new
# "cfop hours" sheet summarizing cfop hours/percentage, mirroring the
# existing "department hours" / "unit(accumulative) hours" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "PI hours" sheet: add a new "cfop" column (G) next to "app" (F).
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("PI hours")

$ws1.Range("G1").Value = "cfop"
$ws1.Range("G2").Value = "['cfop_NH']"
$ws1.Range("G3").Value = "['cfop_CHOUDHURY', 'cfop_RRC']"

# Copy the header formatting (bold, border, centered) from an existing
# header cell onto the new header cell, values only, just paste formats.
$ws1.Range("B1").Copy() | Out-Null
$ws1.Range("G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. New "cfop hours" sheet, appended after the last existing sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "cfop hours"

$ws4.Range("B1").Value = "cfop"
$ws4.Range("C1").Value = "hours"
$ws4.Range("D1").Value = "percentage"

$ws4.Range("A2").Value = 0
$ws4.Range("B2").Value = "cfop_RRC"
$ws4.Range("C2").Value = 74
$ws4.Range("D2").Value = 42.40687679083094

$ws4.Range("A3").Value = 1
$ws4.Range("B3").Value = "cfop_NH"
$ws4.Range("C3").Value = 63.5
$ws4.Range("D3").Value = 36.38968481375358

$ws4.Range("A4").Value = 2
$ws4.Range("B4").Value = "cfop_CHOUDHURY"
$ws4.Range("C4").Value = 37
$ws4.Range("D4").Value = 21.20343839541547

# Copy header formatting (bold, border, centered) + the "A" column
# index formatting from the "unit(accumulative) hours" sheet, which has
# the same B/C/D header + A-index-column layout as the new sheet.
$ws3 = $wb.Worksheets.Item("unit(accumulative) hours")
$ws3.Range("B1:D1").Copy() | Out-Null
$ws4.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws3.Range("A2:A4").Copy() | Out-Null
$ws4.Range("A2:A4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Restore the original active sheet / tab selection.
$ws1.Activate()
